# Update the "想去人数" (want-to-go count) values in column F across all 4 sheets
# to reflect the refreshed snapshot committed in 456a3b4.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 796   # F3: 795 -> 796
$ws.Cells.Item(5, 6).Value = 428   # F5: 427 -> 428
$ws.Cells.Item(6, 6).Value = 738   # F6: 737 -> 738
$ws.Cells.Item(8, 6).Value = 930   # F8: 928 -> 930
$ws.Cells.Item(10, 6).Value = 986   # F10: 984 -> 986
$ws.Cells.Item(14, 6).Value = 98   # F14: 95 -> 98
$ws.Cells.Item(17, 6).Value = 24909   # F17: 24888 -> 24909
$ws.Cells.Item(18, 6).Value = 2443   # F18: 2437 -> 2443
$ws.Cells.Item(19, 6).Value = 159   # F19: 158 -> 159
$ws.Cells.Item(22, 6).Value = 102   # F22: 99 -> 102
$ws.Cells.Item(25, 6).Value = 133   # F25: 128 -> 133
$ws.Cells.Item(26, 6).Value = 256   # F26: 255 -> 256
$ws.Cells.Item(29, 6).Value = 52   # F29: 51 -> 52
$ws.Cells.Item(30, 6).Value = 378   # F30: 376 -> 378
$ws.Cells.Item(31, 6).Value = 27   # F31: 26 -> 27
$ws.Cells.Item(32, 6).Value = 463   # F32: 462 -> 463
$ws.Cells.Item(33, 6).Value = 197   # F33: 196 -> 197
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(7, 6).Value = 246   # F7: 242 -> 246
$ws.Cells.Item(10, 6).Value = 3694   # F10: 3692 -> 3694
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 282   # F2: 281 -> 282
$ws.Cells.Item(4, 6).Value = 839   # F4: 836 -> 839
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 282   # F2: 281 -> 282
$ws.Cells.Item(5, 6).Value = 839   # F5: 836 -> 839
$ws.Cells.Item(6, 6).Value = 796   # F6: 795 -> 796
$ws.Cells.Item(8, 6).Value = 428   # F8: 427 -> 428
$ws.Cells.Item(9, 6).Value = 738   # F9: 737 -> 738
$ws.Cells.Item(16, 6).Value = 930   # F16: 928 -> 930
$ws.Cells.Item(18, 6).Value = 986   # F18: 984 -> 986
$ws.Cells.Item(21, 6).Value = 98   # F21: 95 -> 98
$ws.Cells.Item(24, 6).Value = 24910   # F24: 24888 -> 24910
$ws.Cells.Item(31, 6).Value = 2443   # F31: 2437 -> 2443
$ws.Cells.Item(32, 6).Value = 159   # F32: 158 -> 159
$ws.Cells.Item(39, 6).Value = 256   # F39: 255 -> 256
$ws.Cells.Item(43, 6).Value = 52   # F43: 51 -> 52
$ws.Cells.Item(46, 6).Value = 27   # F46: 26 -> 27
$ws.Cells.Item(47, 6).Value = 463   # F47: 462 -> 463
$ws.Cells.Item(48, 6).Value = 197   # F48: 196 -> 197
